$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2982.7
$ws.Range("I51").Value = 1543
$ws.Range("J51").Value = 3599.7144
$ws.Range("K51").Value = 1543
$ws.Range("L51").Value = 3599.7144
$ws.Range("M51").Value = -1059
$ws.Range("N51").Value = -4567.7144
$ws.Range("H96").Value = 3872.6
$ws.Range("I96").Value = 3303
$ws.Range("K96").Value = 9909
$ws.Range("M96").Value = -8536
$ws.Range("H116").Value = 5059.8823
$ws.Range("I116").Value = 4262.2
$ws.Range("K116").Value = 4262.2
$ws.Range("M116").Value = -820.1999999999998
$ws.Range("H138").Value = 6807734.5
$ws.Range("I138").Value = 2356.6
$ws.Range("J138").Value = 8552703
$ws.Range("K138").Value = 7069.799999999999
$ws.Range("L138").Value = 25658109
$ws.Range("M138").Value = -1929.799999999999
$ws.Range("N138").Value = -25668389

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12354713
$ws.Range("I32").Value = 16952944
$ws.Range("K32").Value = 16952944
$ws.Range("M32").Value = -16952657
$ws.Range("H61").Value = 55560708
$ws.Range("I61").Value = 100003450
$ws.Range("J61").Value = 7275
$ws.Range("K61").Value = 100003450
$ws.Range("L61").Value = 7275
$ws.Range("M61").Value = -100003238
$ws.Range("N61").Value = -7699
$ws.Range("H75").Value = 66666.664
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61748
$ws.Range("H78").Value = 66666.664
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -188736
$ws.Range("H132").Value = 52634850
$ws.Range("I132").Value = 3451.1667
$ws.Range("K132").Value = 10353.5001
$ws.Range("M132").Value = -7823.500100000001
$ws.Range("H136").Value = 55560708
$ws.Range("I136").Value = 100003450
$ws.Range("J136").Value = 7275
$ws.Range("K136").Value = 300010350
$ws.Range("L136").Value = 21825
$ws.Range("M136").Value = -300007800
$ws.Range("N136").Value = -26925
$ws.Range("H139").Value = 39999.5
$ws.Range("J139").Value = 39999.5
$ws.Range("L139").Value = 39999.5
$ws.Range("N139").Value = -50279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 4249.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H33").Value = 5000
$ws.Range("J33").Value = 5000
$ws.Range("L33").Value = 5000
$ws.Range("N33").Value = -5672
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""
$ws.Range("H81").Value = 44776.125
$ws.Range("J81").Value = 44776.125
$ws.Range("L81").Value = 44776.125
$ws.Range("N81").Value = -46898.125
$ws.Range("H84").Value = 44776.125
$ws.Range("J84").Value = 44776.125
$ws.Range("L84").Value = 134328.375
$ws.Range("N84").Value = -144936.375
$ws.Range("H94").Value = 1801.25
$ws.Range("J94").Value = 2005
$ws.Range("L94").Value = 2005
$ws.Range("N94").Value = -2907

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45460852
$ws.Range("I31").Value = 6128.1816
$ws.Range("J31").Value = 90915576
$ws.Range("K31").Value = 6128.1816
$ws.Range("L31").Value = 90915576
$ws.Range("M31").Value = -5833.1816
$ws.Range("N31").Value = -90916166
$ws.Range("H34").Value = 45460852
$ws.Range("I34").Value = 6128.1816
$ws.Range("J34").Value = 90915576
$ws.Range("K34").Value = 6128.1816
$ws.Range("L34").Value = 90915576
$ws.Range("M34").Value = -5926.1816
$ws.Range("N34").Value = -90915980
$ws.Range("H58").Value = 5171
$ws.Range("I58").Value = 4079.8
$ws.Range("K58").Value = 4079.8
$ws.Range("M58").Value = -3876.8
$ws.Range("H86").Value = 3108.2632
$ws.Range("I86").Value = 2581.077
$ws.Range("K86").Value = 2581.077
$ws.Range("M86").Value = -1458.077
$ws.Range("H89").Value = 3108.2632
$ws.Range("I89").Value = 2581.077
$ws.Range("K89").Value = 12905.385
$ws.Range("M89").Value = -7289.385000000002
$ws.Range("H99").Value = 8594.583000000001
$ws.Range("I99").Value = 9063.799999999999
$ws.Range("J99").Value = 6248.5
$ws.Range("K99").Value = 9063.799999999999
$ws.Range("L99").Value = 6248.5
$ws.Range("M99").Value = -7565.799999999999
$ws.Range("N99").Value = -9244.5
$ws.Range("H107").Value = 1250.3158
$ws.Range("I107").Value = 769.0714
$ws.Range("K107").Value = 769.0714
$ws.Range("M107").Value = 1150.9286
$ws.Range("H126").Value = 8594.583000000001
$ws.Range("I126").Value = 9063.799999999999
$ws.Range("J126").Value = 6248.5
$ws.Range("K126").Value = 27191.4
$ws.Range("L126").Value = 18745.5
$ws.Range("M126").Value = -24721.4
$ws.Range("N126").Value = -23685.5
$ws.Range("H136").Value = 5171
$ws.Range("I136").Value = 4079.8
$ws.Range("K136").Value = 12239.4
$ws.Range("M136").Value = -9689.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 50.714287
$ws.Range("I11").Value = 50.714287
$ws.Range("K11").Value = 152.142861
$ws.Range("M11").Value = -12.14286099999998
$ws.Range("H56").Value = 32351.125
$ws.Range("I56").Value = 32351.125
$ws.Range("K56").Value = 32351.125
$ws.Range("M56").Value = -31821.125
$ws.Range("H122").Value = 1812.5
$ws.Range("I122").Value = 868.8
$ws.Range("J122").Value = 2486.5715
$ws.Range("K122").Value = 7819.2
$ws.Range("L122").Value = 22379.1435
$ws.Range("M122").Value = -5369.2
$ws.Range("N122").Value = -27279.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 16673320
$ws.Range("I126").Value = 14295652
$ws.Range("J126").Value = 18186382
$ws.Range("K126").Value = 42886956
$ws.Range("L126").Value = 54559146
$ws.Range("M126").Value = -42884486
$ws.Range("N126").Value = -54564086
$ws.Range("H132").Value = 2286.111
$ws.Range("I132").Value = 2422.3333
$ws.Range("J132").Value = 2013.6666
$ws.Range("K132").Value = 7266.999899999999
$ws.Range("L132").Value = 6040.9998
$ws.Range("M132").Value = -4736.999899999999
$ws.Range("N132").Value = -11100.9998
$ws.Range("H133").Value = 142723.1
$ws.Range("J133").Value = 142723.1
$ws.Range("L133").Value = 142723.1
$ws.Range("N133").Value = -152843.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5613.857
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 5882.8335
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 5882.8335
$ws.Range("M22").Value = -3705
$ws.Range("N22").Value = -6472.8335
$ws.Range("H27").Value = 5613.857
$ws.Range("I27").Value = 4000
$ws.Range("J27").Value = 5882.8335
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 5882.8335
$ws.Range("M27").Value = -3893
$ws.Range("N27").Value = -6096.8335
$ws.Range("H40").Value = 5010.393
$ws.Range("I40").Value = 4489.1816
$ws.Range("J40").Value = 6921.5
$ws.Range("K40").Value = 4489.1816
$ws.Range("L40").Value = 6921.5
$ws.Range("M40").Value = -4353.1816
$ws.Range("N40").Value = -7193.5
$ws.Range("H55").Value = 1290.5454
$ws.Range("J55").Value = 1708.1666
$ws.Range("L55").Value = 1708.1666
$ws.Range("N55").Value = -2054.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5453.4443
$ws.Range("I96").Value = 2959.5
$ws.Range("K96").Value = 2959.5
$ws.Range("M96").Value = -1586.5
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = ""
$ws.Range("N101").Value = ""
$ws.Range("H102").Value = 65000
$ws.Range("I102").Value = 60000
$ws.Range("J102").Value = 70000
$ws.Range("K102").Value = 60000
$ws.Range("L102").Value = 70000
$ws.Range("M102").Value = -56755
$ws.Range("N102").Value = -76490
$ws.Range("H103").Value = 68294.8
$ws.Range("J103").Value = 68294.8
$ws.Range("L103").Value = 68294.8
$ws.Range("N103").Value = -70638.8
$ws.Range("H106").Value = 37571.832
$ws.Range("I106").Value = 32000
$ws.Range("J106").Value = 38686.2
$ws.Range("K106").Value = 32000
$ws.Range("L106").Value = 38686.2
$ws.Range("M106").Value = -30738
$ws.Range("N106").Value = -41210.2
$ws.Range("H107").Value = 513.38464
$ws.Range("I107").Value = 489.9
$ws.Range("K107").Value = 1469.7
$ws.Range("M107").Value = 450.3000000000002
$ws.Range("H126").Value = 4220.185
$ws.Range("I126").Value = 5137.25
$ws.Range("K126").Value = 15411.75
$ws.Range("M126").Value = -12941.75
$ws.Range("H132").Value = 4159.3228
$ws.Range("I132").Value = 4320.345
$ws.Range("J132").Value = 2894.1428
$ws.Range("K132").Value = 12961.035
$ws.Range("L132").Value = 8682.428400000001
$ws.Range("M132").Value = -10431.035
$ws.Range("N132").Value = -13742.4284
$ws.Range("H136").Value = 1586.2325
$ws.Range("I136").Value = 1576.381
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4729.143
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2179.143
$ws.Range("N136").Value = -11100
